# Updates the lattice-multiplication practice table: replaces the
# multiplication problem (and its partial products) in each of the
# 15 table cells (5 rows x 3 columns) with new values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11   # soft line break (w:br) inside a Word range

function Set-Problem {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Top,
        [string]$Side,
        [string]$Digit1,
        [string]$Digit2
    )

    $cell = $t.Cell($Row, $Col)
    $text = $Top + $vt + $Side + $vt + "  ----" + $vt + $Digit1 + $vt + $Digit2
    $cell.Range.Text = $text
}

# Row 1
Set-Problem 1 1 "36 x 49" "  4    9" "3|    |" "6|    |"
Set-Problem 1 2 "63 x 64" "  6    4" "6|    |" "3|    |"
Set-Problem 1 3 "10 x 74" "  7    4" "1|    |" "0|    |"

# Row 2
Set-Problem 2 1 "86 x 64" "  6    4" "8|    |" "6|    |"
Set-Problem 2 2 "99 x 45" "  4    5" "9|    |" "9|    |"
Set-Problem 2 3 "60 x 32" "  3    2" "6|    |" "0|    |"

# Row 3
Set-Problem 3 1 "68 x 39" "  3    9" "6|    |" "8|    |"
Set-Problem 3 2 "71 x 54" "  5    4" "7|    |" "1|    |"
Set-Problem 3 3 "20 x 28" "  2    8" "2|    |" "0|    |"

# Row 4
Set-Problem 4 1 "23 x 19" "  1    9" "2|    |" "3|    |"
Set-Problem 4 2 "13 x 71" "  7    1" "1|    |" "3|    |"
Set-Problem 4 3 "88 x 83" "  8    3" "8|    |" "8|    |"

# Row 5
Set-Problem 5 1 "83 x 90" "  9    0" "8|    |" "3|    |"
Set-Problem 5 2 "48 x 94" "  9    4" "4|    |" "8|    |"
Set-Problem 5 3 "85 x 65" "  6    5" "8|    |" "5|    |"
